$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helper: find the first shape in a Shapes collection whose PlaceholderFormat
# is the "Date" placeholder (ppPlaceholderDate = 16). Some shapes in the
# collection are not placeholders at all, so guard with try/catch.
# ---------------------------------------------------------------------------
function Get-DatePlaceholderShape($shapeCollection) {
    for ($i = 1; $i -le $shapeCollection.Count; $i++) {
        $shp = $shapeCollection.Item($i)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePh = $true
            }
        } catch {
        }
        if ($isDatePh) {
            return $shp
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Helper: replace one exact occurrence of $oldText with $newText inside a
# shape's TextFrame, using a Characters() sub-range that spans exactly the
# matched text. Using a range that aligns with the full run keeps the
# surrounding run(s)/field intact instead of fragmenting them.
# ---------------------------------------------------------------------------
function Set-ExactTextRange($shape, [string]$oldText, [string]$newText) {
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        Write-Host "WARNING: text not found -> " $oldText
        return
    }
    $sub = $tr.Characters($idx + 1, $oldText.Length)
    $sub.Text = $newText
}

$oldDate = "6/5/2022"
$newDate = "12/27/2022"

# 1) Update the fixed "datetimeFigureOut" date field text everywhere it is
#    defined: once on the slide master, and once on each of the 11 custom
#    (slide) layouts.
$masterDateShape = Get-DatePlaceholderShape $p.SlideMaster.Shapes
if ($masterDateShape -ne $null) {
    Set-ExactTextRange $masterDateShape $oldDate $newDate
}

for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    $layoutDateShape = Get-DatePlaceholderShape $layout.Shapes
    if ($layoutDateShape -ne $null) {
        Set-ExactTextRange $layoutDateShape $oldDate $newDate
    }
}

# 2) Slide 10 ("Object of Analysis"): "Doctors" -> "physicians"
$slide10 = $p.Slides.Item(10)
$slide10Body = $slide10.Shapes.Item(2)
Set-ExactTextRange $slide10Body `
    "The objective of this analysis is to provide a model(s) to assist Doctors in diagnosing heart disease based on the features present.  " `
    "The objective of this analysis is to provide a model(s) to assist physicians in diagnosing heart disease based on the features present.  "

# 3) Slide 5 ("Feature Descriptions - Clarifications"): capitalize "brought"
$slide5 = $p.Slides.Item(5)
$slide5Body = $slide5.Shapes.Item(2)
Set-ExactTextRange $slide5Body `
    "brought on by exertion or emotional stress" `
    "Brought on by exertion or emotional stress"
